$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 247
$ws.Range("I107").Value = 198.80952
$ws.Range("J107").Value = 753
$ws.Range("K107").Value = 198.80952
$ws.Range("L107").Value = 753
$ws.Range("M107").Value = 1721.19048
$ws.Range("N107").Value = -4593
$ws.Range("H112").Value = 941.7778
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 941.7778
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 2825.3334
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -5041.3334
$ws.Range("H137").Value = 2133.8462
$ws.Range("I137").Value = 1432.3889
$ws.Range("J137").Value = 3712.125
$ws.Range("K137").Value = 4297.1667
$ws.Range("L137").Value = 11136.375
$ws.Range("M137").Value = -1747.1667
$ws.Range("N137").Value = -16236.375
$ws.Range("H138").Value = 6618.646
$ws.Range("I138").Value = 1297.0857
$ws.Range("J138").Value = 20945.924
$ws.Range("K138").Value = 3891.2571
$ws.Range("L138").Value = 62837.772
$ws.Range("M138").Value = 1248.7429
$ws.Range("N138").Value = -73117.772

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 41778.37
$ws.Range("I32").Value = 7104.1763
$ws.Range("K32").Value = 7104.1763
$ws.Range("M32").Value = -6817.1763
$ws.Range("H61").Value = 2170.611
$ws.Range("I61").Value = 1508.1428
$ws.Range("J61").Value = 2592.182
$ws.Range("K61").Value = 1508.1428
$ws.Range("L61").Value = 2592.182
$ws.Range("M61").Value = -1296.1428
$ws.Range("N61").Value = -3016.182
$ws.Range("H74").Value = 1199.3667
$ws.Range("I74").Value = 1060.8667
$ws.Range("J74").Value = 1337.8667
$ws.Range("K74").Value = 1060.8667
$ws.Range("L74").Value = 1337.8667
$ws.Range("M74").Value = -186.8667
$ws.Range("N74").Value = -3085.8667
$ws.Range("H77").Value = 1199.3667
$ws.Range("I77").Value = 1060.8667
$ws.Range("J77").Value = 1337.8667
$ws.Range("K77").Value = 5304.333500000001
$ws.Range("L77").Value = 6689.333500000001
$ws.Range("M77").Value = -936.3335000000006
$ws.Range("N77").Value = -15425.3335
$ws.Range("H122").Value = 1248.9166
$ws.Range("I122").Value = 1174.85
$ws.Range("K122").Value = 3524.55
$ws.Range("M122").Value = -1074.55
$ws.Range("H136").Value = 2170.611
$ws.Range("I136").Value = 1508.1428
$ws.Range("J136").Value = 2592.182
$ws.Range("K136").Value = 4524.428400000001
$ws.Range("L136").Value = 7776.545999999999
$ws.Range("M136").Value = -1974.428400000001
$ws.Range("N136").Value = -12876.546

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 45088.31
$ws.Range("I86").Value = 93724.586
$ws.Range("J86").Value = 3400.0715
$ws.Range("K86").Value = 93724.586
$ws.Range("L86").Value = 3400.0715
$ws.Range("M86").Value = -92601.586
$ws.Range("N86").Value = -5646.0715
$ws.Range("H89").Value = 45088.31
$ws.Range("I89").Value = 93724.586
$ws.Range("J89").Value = 3400.0715
$ws.Range("K89").Value = 468622.93
$ws.Range("L89").Value = 17000.3575
$ws.Range("M89").Value = -463006.93
$ws.Range("N89").Value = -28232.3575
$ws.Range("H137").Value = 40775
$ws.Range("J137").Value = 40775
$ws.Range("L137").Value = 40775
$ws.Range("N137").Value = -50975

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30213.781
$ws.Range("I31").Value = 988
$ws.Range("J31").Value = 49697.637
$ws.Range("K31").Value = 988
$ws.Range("L31").Value = 49697.637
$ws.Range("M31").Value = -693
$ws.Range("N31").Value = -50287.637
$ws.Range("H34").Value = 30213.781
$ws.Range("I34").Value = 988
$ws.Range("J34").Value = 49697.637
$ws.Range("K34").Value = 988
$ws.Range("L34").Value = 49697.637
$ws.Range("M34").Value = -786
$ws.Range("N34").Value = -50101.637
$ws.Range("H99").Value = 13338.454
$ws.Range("I99").Value = 4340
$ws.Range("J99").Value = 29085.75
$ws.Range("K99").Value = 4340
$ws.Range("L99").Value = 29085.75
$ws.Range("M99").Value = -2842
$ws.Range("N99").Value = -32081.75
$ws.Range("H126").Value = 13338.454
$ws.Range("I126").Value = 4340
$ws.Range("J126").Value = 29085.75
$ws.Range("K126").Value = 13020
$ws.Range("L126").Value = 87257.25
$ws.Range("M126").Value = -10550
$ws.Range("N126").Value = -92197.25
$ws.Range("H131").Value = 14058.517
$ws.Range("J131").Value = 14058.517
$ws.Range("L131").Value = 14058.517
$ws.Range("N131").Value = -24138.517

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 632.5
$ws.Range("I34").Value = 120
$ws.Range("J34").Value = 940
$ws.Range("K34").Value = 360
$ws.Range("L34").Value = 2820
$ws.Range("M34").Value = -276
$ws.Range("N34").Value = -2988
$ws.Range("H44").Value = 1024.6522
$ws.Range("I44").Value = 299.8
$ws.Range("J44").Value = 1226
$ws.Range("K44").Value = 899.4000000000001
$ws.Range("L44").Value = 3678
$ws.Range("M44").Value = -501.4000000000001
$ws.Range("N44").Value = -4474

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("H26").Value = 7557.727
$ws.Range("J26").Value = 7557.727
$ws.Range("L26").Value = 7557.727
$ws.Range("N26").Value = -8117.727
$ws.Range("H50").Value = 7557.727
$ws.Range("J50").Value = 7557.727
$ws.Range("L50").Value = 7557.727
$ws.Range("N50").Value = -8553.726999999999
$ws.Range("H53").Value = 12405.714
$ws.Range("I53").Value = 10520
$ws.Range("J53").Value = 13160
$ws.Range("K53").Value = 10520
$ws.Range("L53").Value = 13160
$ws.Range("M53").Value = -9889
$ws.Range("N53").Value = -14422
$ws.Range("H102").Value = 2579.9092
$ws.Range("I102").Value = 2114.1667
$ws.Range("J102").Value = 3138.8
$ws.Range("K102").Value = 2114.1667
$ws.Range("L102").Value = 3138.8
$ws.Range("M102").Value = -492.1667000000002
$ws.Range("N102").Value = -6382.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5800
$ws.Range("I7").Value = 2400
$ws.Range("K7").Value = 2400
$ws.Range("M7").Value = -2288
$ws.Range("H68").Value = 2909.1875
$ws.Range("I68").Value = 1457.1428
$ws.Range("J68").Value = 4038.5557
$ws.Range("K68").Value = 1457.1428
$ws.Range("L68").Value = 4038.5557
$ws.Range("M68").Value = -708.1428000000001
$ws.Range("N68").Value = -5536.5557
$ws.Range("H71").Value = 2909.1875
$ws.Range("I71").Value = 1457.1428
$ws.Range("J71").Value = 4038.5557
$ws.Range("K71").Value = 7285.714
$ws.Range("L71").Value = 20192.7785
$ws.Range("M71").Value = -3541.714
$ws.Range("N71").Value = -27680.7785
$ws.Range("H126").Value = 5800
$ws.Range("I126").Value = 2400
$ws.Range("K126").Value = 7200
$ws.Range("M126").Value = -4730
$ws.Range("H136").Value = 1936.9412
$ws.Range("I136").Value = 1509.3636
$ws.Range("J136").Value = 2720.8333
$ws.Range("K136").Value = 4528.0908
$ws.Range("L136").Value = 8162.499899999999
$ws.Range("M136").Value = -1978.0908
$ws.Range("N136").Value = -13262.4999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2302.8333
$ws.Range("I126").Value = 2246
$ws.Range("J126").Value = 2416.5
$ws.Range("K126").Value = 6738
$ws.Range("L126").Value = 7249.5
$ws.Range("M126").Value = -4268
$ws.Range("N126").Value = -12189.5
$ws.Range("H132").Value = 2420.0232
$ws.Range("I132").Value = 2227.075
$ws.Range("J132").Value = 4992.6665
$ws.Range("K132").Value = 6681.224999999999
$ws.Range("L132").Value = 14977.9995
$ws.Range("M132").Value = -4151.224999999999
$ws.Range("N132").Value = -20037.9995
$ws.Range("H136").Value = 795.8421
$ws.Range("I136").Value = 464.26666
$ws.Range("J136").Value = 2039.25
$ws.Range("K136").Value = 1392.79998
$ws.Range("L136").Value = 6117.75
$ws.Range("M136").Value = 1157.20002
$ws.Range("N136").Value = -11217.75
